$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lookup Functions")

# --- Remove the blank row above the Xlookup table, shifting rows 49-57 up to 48-56 ---
$ws.Rows("48").Delete() | Out-Null

# --- Workbook window position ---
$win = $excel.ActiveWindow
$win.Left = 17520
$win.Top = 3120

# --- Sheet view: hide gridlines, zoom to 50%, change selection ---
$win.DisplayGridlines = $false
$win.Zoom = 50
$ws.Range("K21").Select() | Out-Null

# --- New column width for column I ---
$ws.Columns.Item(9).ColumnWidth = 12

# --- Print area ---
$ws.PageSetup.PrintArea = '$B$2:$I$56'

# --- Page setup: scale + orientation + headings/gridlines for print ---
$ps = $ws.PageSetup
$ps.Zoom = 60
$ps.Orientation = 1
$ps.PrintHeadings = $true
$ps.PrintGridlines = $true

# --- Manual page breaks ---
$ws.Rows.Item(2).PageBreak = 1
$ws.Rows.Item(57).PageBreak = 1
$ws.Columns.Item(10).PageBreak = 1
